# Update NATMI LR-pair stats for Dlk1-Notch2 with refreshed TPM-based numbers.
# Shared-string table content: new row set references ECs/FAPs/MuSCs/Resolving-Mac
# as both sending and target clusters (previously Resolving-Mac only appeared as a
# target cluster; now the ligand/receptor pair is also reported with Resolving-Mac
# as the sending cluster, rows 14-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E-J ligand stats, K-P receptor stats, Q-T edge stats.
$rows = @(
    @("ECs", "Dlk1", "Notch2", "ECs", 3, 1, 0.506144, 1.518432, 0.002728503728229319, 0.002728503728229319, 3, 1, 1.174933333333333, 3.5248, 0.01171850713626266, 0.01171850713626266, 0.5946854570666668, 5.3521691136, ([double]"3.197399041057454e-05"), ([double]"3.197399041057453e-05")),
    @("ECs", "Dlk1", "Notch2", "FAPs", 3, 1, 0.506144, 1.518432, 0.002728503728229319, 0.002728503728229319, 3, 1, 27.50472933333333, 82.514188, 0.2743256641287217, 0.2743256641287218, 13.92135372369067, 125.292183513216, 0.0007484985973242011, 0.0007484985973242013),
    @("ECs", "Dlk1", "Notch2", "MuSCs", 3, 1, 0.506144, 1.518432, 0.002728503728229319, 0.002728503728229319, 3, 1, 39.361408, 118.084224, 0.3925813724534833, 0.3925813724534833, 19.922540490752, 179.302864416768, 0.001071159738372712, 0.001071159738372712),
    @("ECs", "Dlk1", "Notch2", "Resolving-Mac", 3, 1, 0.506144, 1.518432, 0.002728503728229319, 0.002728503728229319, 3, 1, 32.221985, 96.665955, 0.3213744562815322, 0.3213744562815322, 16.30896437584, 146.78067938256, 0.0008768714021218306, 0.0008768714021218308),
    @("FAPs", "Dlk1", "Notch2", "ECs", 3, 1, 174.184255, 522.552765, 0.9389865120723476, 0.9389865120723477, 3, 1, 1.174933333333333, 3.5248, 0.01171850713626266, 0.01171850713626266, 204.6548873413334, 1841.893986072, 0.01100352014257419, 0.01100352014257419),
    @("FAPs", "Dlk1", "Notch2", "FAPs", 3, 1, 174.184255, 522.552765, 0.9389865120723476, 0.9389865120723477, 3, 1, 27.50472933333333, 82.514188, 0.2743256641287217, 0.2743256641287218, 4790.890787903313, 43118.01709112983, 0.2575880985321587, 0.2575880985321588),
    @("FAPs", "Dlk1", "Notch2", "MuSCs", 3, 1, 174.184255, 522.552765, 0.9389865120723476, 0.9389865120723477, 3, 1, 39.361408, 118.084224, 0.3925813724534833, 0.3925813724534833, 6856.137528231041, 61705.23775407937, 0.3686286136246715, 0.3686286136246716),
    @("FAPs", "Dlk1", "Notch2", "Resolving-Mac", 3, 1, 174.184255, 522.552765, 0.9389865120723476, 0.9389865120723477, 3, 1, 32.221985, 96.665955, 0.3213744562815322, 0.3213744562815322, 5612.562451846175, 50513.06206661557, 0.3017662797729431, 0.3017662797729431),
    @("MuSCs", "Dlk1", "Notch2", "ECs", 3, 1, 10.788518, 32.365554, 0.05815837308171015, 0.05815837308171016, 3, 1, 1.174933333333333, 3.5248, 0.01171850713626266, 0.01171850713626266, 12.67578941546667, 114.0821047392, 0.0006815293099914464, 0.0006815293099914464),
    @("MuSCs", "Dlk1", "Notch2", "FAPs", 3, 1, 10.788518, 32.365554, 0.05815837308171015, 0.05815837308171016, 3, 1, 27.50472933333333, 82.514188, 0.2743256641287217, 0.2743256641287218, 296.7352674977946, 2670.617407480152, 0.01595433432028611, 0.01595433432028611),
    @("MuSCs", "Dlk1", "Notch2", "MuSCs", 3, 1, 10.788518, 32.365554, 0.05815837308171015, 0.05815837308171016, 3, 1, 39.361408, 118.084224, 0.3925813724534833, 0.3925813724534833, 424.651258713344, 3821.861328420096, 0.02283189392407949, 0.0228318939240795),
    @("MuSCs", "Dlk1", "Notch2", "Resolving-Mac", 3, 1, 10.788518, 32.365554, 0.05815837308171015, 0.05815837308171016, 3, 1, 32.221985, 96.665955, 0.3213744562815322, 0.3213744562815322, 347.6274651682299, 3128.647186514069, 0.0186906155273531, 0.0186906155273531),
    @("Resolving-Mac", "Dlk1", "Notch2", "ECs", 1, 0.3333333333333333, 0.02348666666666667, 0.07046, 0.0001266111177129024, 0.0001266111177129024, 3, 1, 1.174933333333333, 3.5248, 0.01171850713626266, 0.01171850713626266, 0.02759526755555556, 0.248357408, ([double]"1.483693286448838e-06"), ([double]"1.483693286448837e-06")),
    @("Resolving-Mac", "Dlk1", "Notch2", "FAPs", 1, 0.3333333333333333, 0.02348666666666667, 0.07046, 0.0001266111177129024, 0.0001266111177129024, 3, 1, 27.50472933333333, 82.514188, 0.2743256641287217, 0.2743256641287218, 0.6459944096088889, 5.81394968648, ([double]"3.47326789526717e-05"), ([double]"3.473267895267171e-05")),
    @("Resolving-Mac", "Dlk1", "Notch2", "MuSCs", 1, 0.3333333333333333, 0.02348666666666667, 0.07046, 0.0001266111177129024, 0.0001266111177129024, 3, 1, 39.361408, 118.084224, 0.3925813724534833, 0.3925813724534833, 0.9244682692266667, 8.32021442304, ([double]"4.970516635960074e-05"), ([double]"4.970516635960074e-05")),
    @("Resolving-Mac", "Dlk1", "Notch2", "Resolving-Mac", 1, 0.3333333333333333, 0.02348666666666667, 0.07046, 0.0001266111177129024, 0.0001266111177129024, 3, 1, 32.221985, 96.665955, 0.3213744562815322, 0.3213744562815322, 0.7567870210333333, 6.8110831893, ([double]"4.068957911418106e-05"), ([double]"4.068957911418106e-05"))
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value2 = $values[$j]
    }
}
